$d = $word.ActiveDocument

# Locate the paragraph that ends the "Commit 19" description (contains
# ", hasOwnProperty()") so the new "Commit 20" content can be inserted
# right after it, before the following (empty, numbered) paragraph.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*hasOwnProperty()*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find anchor paragraph containing 'hasOwnProperty()'"
}

$anchorPara = $d.Paragraphs.Item($targetIndex)
$r = $anchorPara.Range
$r.Collapse(0)

# New paragraph: "Commit 20:"
$r.InsertParagraphAfter()
$r.Collapse(0)
$r.InsertAfter("Commit 20:")

# New paragraph with the commit's description, split across two runs (as
# it was authored in two separate edits in the original document).
$r.Collapse(0)
$r.InsertParagraphAfter()
$r.Collapse(0)
$null = $r.Move(1, 1)
$run1Start = $r.Start
$r.InsertAfter("Prototype in array, Object.getPrototypeOf(), Array.prototype, changing function prototype from object to array")
$run1End = $r.End
$r.Collapse(0)
$r.InsertAfter(", classes, extends and inheritance")

# Keep the two pieces of text as separate runs (rather than letting the
# engine silently merge them because they share identical formatting) by
# toggling a character-formatting property on the first run back to its
# original value - this has no visible effect but prevents run-coalescing.
$run1 = $d.Range($run1Start, $run1End)
$run1.Font.Bold = 1
$run1.Font.Bold = 0
